$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.344.40'
$ws.Range("E2").Value = '  -0.04%  '

$ws.Range("D3").Value = '2.327.04'
$ws.Range("E3").Value = '  -1.21%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").Value = '304.16'
$ws.Range("E5").Value = '  -2.04%  '

$ws.Range("D6").Value = '100.94'
$ws.Range("E6").Value = '  -2.92%  '

$ws.Range("D7").Value = '0.509'
$ws.Range("E7").Value = '  -3.38%  '

$ws.Range("D9").Value = '0.506'
$ws.Range("E9").Value = '  -3.27%  '

$ws.Range("D10").Value = '35.18'
$ws.Range("E10").Value = '  -2.69%  '

$ws.Range("D11").Value = '0.0797'
$ws.Range("E11").Value = '  -2.08%  '

$ws.Range("E12").Value = '  +0.41%  '

$ws.Range("D13").Value = '6.78'
$ws.Range("E13").Value = '  -3.17%  '

$ws.Range("D14").Value = '2.695.45'
$ws.Range("E14").Value = '  -1.09%  '

$ws.Range("D15").Value = '15.59'
$ws.Range("E15").Value = '  -0.89%  '

$ws.Range("D16").Value = '2.340.74'
$ws.Range("E16").Value = '  -1.41%  '

$ws.Range("D17").Value = '0.803'
$ws.Range("E17").Value = '  -1.21%  '

$ws.Range("D18").Value = '43.283.69'
$ws.Range("E18").Value = '  -0.07%  '

$ws.Range("D19").Value = '11.81'
$ws.Range("E19").Value = '  -1.67%  '

$ws.Range("D20").Value = '0.0₃0908'
$ws.Range("E20").Value = '  -2.23%  '

$ws.Range("D21").Value = '6.10'
$ws.Range("E21").Value = '  -2.69%  '

$ws.Range("D22").Value = '68.08'
$ws.Range("E22").Value = '  -0.34%  '

$ws.Range("D23").Value = '237.31'
$ws.Range("E23").Value = '  -1.94%  '

$ws.Range("D24").Value = '1.97'
$ws.Range("E24").Value = '  -3.98%  '

$ws.Range("D25").Value = '2.53'
$ws.Range("E25").Value = '  -3.43%  '

$ws.Range("E26").Value = '  +0.01%  '

$ws.Range("D27").Value = '24.90'
$ws.Range("E27").Value = '  -4.02%  '

$ws.Range("D28").Value = '2.18'
$ws.Range("E28").Value = '  -5.79%  '

$ws.Range("D29").Value = '34.47'
$ws.Range("E29").Value = '  -5.93%  '

$ws.Range("B30").Value = 'Monero'
$ws.Range("C30").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D30").Value = '165.39'
$ws.Range("E30").Value = '  +1.66%  '

$ws.Range("B31").Value = 'Cosmos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D31").Value = '9.21'
$ws.Range("E31").Value = '  -4.15%  '

$ws.Range("E32").Value = '  +0.11%  '

$ws.Range("D33").Value = '5.05'
$ws.Range("E33").Value = '  -4.59%  '

$ws.Range("D34").Value = '4.54'
$ws.Range("E34").Value = '  -2.11%  '

$ws.Range("E35").Value = '  -4.78%  '

$ws.Range("D36").Value = '16.86'
$ws.Range("E36").Value = '  -7.89%  '

$ws.Range("D37").Value = '0.0704'
$ws.Range("E37").Value = '  -4.96%  '

$ws.Range("D38").Value = '2.90'
$ws.Range("E38").Value = '  -6.79%  '

$ws.Range("D39").Value = '1.81'
$ws.Range("E39").Value = '  -6.37%  '

$ws.Range("D40").Value = '0.103'
$ws.Range("E40").Value = '  -3.75%  '

$ws.Range("E41").Value = '  -3.42%  '

$ws.Range("D42").Value = '2.42'
$ws.Range("E42").Value = '  -0.19%  '

$ws.Range("D43").Value = '1.972.87'
$ws.Range("E43").Value = '  -0.93%  '

$ws.Range("E44").Value = '  -3.03%  '

$ws.Range("D45").Value = '18.56'
$ws.Range("E45").Value = '  -6.19%  '

$ws.Range("D46").Value = '10.02'
$ws.Range("E46").Value = '  -4.03%  '

$ws.Range("D47").Value = '2.92'
$ws.Range("E47").Value = '  -5.74%  '

$ws.Range("D48").Value = '55.80'
$ws.Range("E48").Value = '  -4.59%  '

$ws.Range("D49").Value = '4.80'
$ws.Range("E49").Value = '  +2.29%  '

$ws.Range("D50").Value = '2.554.37'
$ws.Range("E50").Value = '  +0.03%  '

$ws.Range("D51").Value = '1.55'
$ws.Range("E51").Value = '  -2.30%  '
